# Update attendance report: Y2_B2526_Blood_&_lymphatics_session_analysis.xlsx
#
# Applies the session-analysis refresh: a couple of sessions that were
# "Not Recorded" became "Recorded" (1 of 217 students logged), a handful of
# sessions that were "Pending" became "Not Recorded" (their date passed
# without being recorded), three date typos got corrected (day/month
# swapped), and the summary statistics block / per-group breakdown numbers
# were refreshed to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Colors (BGR order, as used by the Excel object model's Interior.Color / Font.Color)
$GREEN  = 9498256   # 90EE90 - "Recorded"
$PINK   = 12695295  # FFB6C1 - "Not Recorded"
$YELLOW = 14745599  # FFFFE0 - "Pending" (unused here, reference only)
$BLACK  = 0

# ---------------------------------------------------------------------
# Class Statistics summary block (K6:L10)
# ---------------------------------------------------------------------
$ws.Range("L6").Value = 2        # Recorded Sessions
$ws.Range("L7").Value = 38       # Missing Sessions
$ws.Range("L8").Value = 113      # Pending Sessions
$ws.Range("L9").Value = "1.3%"   # Coverage %
$ws.Range("L10").Value = "0.5%"  # Average Attendance %

# ---------------------------------------------------------------------
# Helper: mark a session row (columns A:I) as "Recorded" (green) or
# "Not Recorded" (pink), updating the Status cell text and the
# attendance-count cell text, while preserving the existing font.
# ---------------------------------------------------------------------
function Set-SessionRow {
    param(
        [int]$Row,
        [string]$StudentsText,
        [string]$StatusText,
        [int]$Color
    )
    $rowRange = $ws.Range("A" + $Row + ":I" + $Row)
    $rowRange.Interior.Color = $Color
    $rowRange.Font.Color = $BLACK

    if ($StudentsText -ne $null) {
        $ws.Range("H" + $Row).Value = $StudentsText
    }
    $ws.Range("I" + $Row).Value = $StatusText
}

# ---------------------------------------------------------------------
# Row 15 (A1 / PHYSIOLOGY / session 1): Not Recorded -> Recorded
# ---------------------------------------------------------------------
Set-SessionRow -Row 15 -StudentsText "1/217" -StatusText "Recorded" -Color $GREEN
$ws.Range("O15").Value = 2
$ws.Range("P15").Value = 2
$ws.Range("R15").Value = "11.8%"
$ws.Range("S15").Value = "0.5%"

# ---------------------------------------------------------------------
# Row 17 (A1 / POS / session 1): Not Recorded -> Recorded
# ---------------------------------------------------------------------
Set-SessionRow -Row 17 -StudentsText "1/217" -StatusText "Recorded" -Color $GREEN
$ws.Range("P17").Value = 4
$ws.Range("Q17").Value = 13

# ---------------------------------------------------------------------
# Rows 18-23: per-group breakdown counts refreshed (no status change)
# ---------------------------------------------------------------------
$ws.Range("P18").Value = 4
$ws.Range("Q18").Value = 13

$ws.Range("P19").Value = 5
$ws.Range("Q19").Value = 12

$ws.Range("P20").Value = 4
$ws.Range("Q20").Value = 13

$ws.Range("P21").Value = 5
$ws.Range("Q21").Value = 12

$ws.Range("P22").Value = 4
$ws.Range("Q22").Value = 13

$ws.Range("P23").Value = 6
$ws.Range("Q23").Value = 11

# ---------------------------------------------------------------------
# Sessions whose date passed and are now "Not Recorded" (were "Pending")
# ---------------------------------------------------------------------
$pendingToNotRecordedRows = @(51, 68, 76, 81, 83, 100, 105, 109, 122, 126, 148, 154)
foreach ($r in $pendingToNotRecordedRows) {
    Set-SessionRow -Row $r -StudentsText $null -StatusText "Not Recorded" -Color $PINK
}

# ---------------------------------------------------------------------
# Date corrections (day/month had been transposed)
# ---------------------------------------------------------------------
$ws.Range("E76").Value = "10/07/2025"
$ws.Range("E78").Value = "10/12/2025"
$ws.Range("E154").Value = "10/05/2025"
